$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(1010751, 'DGES', 7, '3er. CICLO EBI', 'desconocido', 'F', 'Docente', 'TALLER ARTE - ED.MUSICAL', 41, 136),
    @(1010774, 'DGES', 9, '3er. CICLO EBI', 'desconocido', 'F', 'Docente', 'CS.COMPUTACIONALES (Program.)', 49, 100),
    @(1011136, 'DGES', 8, '3er. CICLO EBI', 'desconocido', 'M', 'Docente', 'HISTORIA', 37, 97),
    @(1011146, 'DGES', 9, '3er. CICLO EBI', 'desconocido', 'F', 'Docente', 'TALLER PARTICIPACION JUVENIL', 20, 59),
    @(1011147, 'DGES', 9, '3er. CICLO EBI', 'desconocido', 'M', 'Docente', 'TALLER PARTICIPACION JUVENIL', 41, 206),
    @(1011183, 'DGES', 9, '3er. CICLO EBI', 'desconocido', 'F', 'Docente', 'TALLER EXPRESION ARTISTICA', 33, 144),
    @(1011190, 'DGES', 8, '3er. CICLO EBI', 'desconocido', 'F', 'Docente', 'ARTE - EDUCACION MUSICAL', 27, 25),
    @(1011196, 'DGES', 8, '3er. CICLO EBI', 'desconocido', 'F', 'Docente', 'CS. FISICO-QUIMICA', 61, 75),
    @(1011198, 'DGES', 9, '3er. CICLO EBI', 'desconocido', 'F', 'Docente', 'INGLES', 63, 102),
    @(1011204, 'DGES', 8, '3er. CICLO EBI', 'desconocido', 'F', 'Docente', 'ARTE - EDUCACION MUSICAL', 55, 207),
    @(1011642, 'DGES', 8, '3er. CICLO EBI', 'desconocido', 'F', 'Docente', 'CS.COMPUTACIONALES (Tecnolog.)', 50, 61),
    @(1011644, 'DGES', 9, '3er. CICLO EBI', 'desconocido', 'M', 'Docente', 'CS.COMPUTACIONALES (Program.)', 42, 56),
    @(1011648, 'DGES', 8, '3er. CICLO EBI', 'desconocido', 'F', 'Docente', 'HISTORIA', 47, 145),
    @(1011717, 'DGES', 7, '3er. CICLO EBI', 'desconocido', 'M', 'Docente', 'LENGUA ESPAÑOLA', 54, 152)
)

$startRow = 15
for ($idx = 0; $idx -lt $data.Count; $idx++) {
    $r = $startRow + $idx
    $row = $data[$idx]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $ws.Cells.Item($r, 9).Value = $row[8]
    $ws.Cells.Item($r, 10).Value = $row[9]
}
